$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "chaymaaa"
$ws.Range("C20").Value = "chachahassayoun@gmail.com"
$ws.Range("D20").Value = "hassayoune"
$ws.Range("E20").Value = 54224709
